$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value2 = 0.2147239263803681
$ws.Range("C2").Value2 = 0.5184049079754601
$ws.Range("J2").Value2 = 0.03067484662576687
$ws.Range("P2").Value2 = 0.1503067484662577
$ws.Range("S2").Value2 = 0.08588957055214724
$ws.Range("B3").Value2 = 0.01149425287356322
$ws.Range("C3").Value2 = 0.02873563218390805
$ws.Range("J3").Value2 = 0.02298850574712644
$ws.Range("P3").Value2 = 0.7471264367816092
$ws.Range("S3").Value2 = 0.1896551724137931
$ws.Range("J4").Value2 = 0.02040816326530612
$ws.Range("P4").Value2 = 0.6326530612244898
$ws.Range("S4").Value2 = 0.3469387755102041
$ws.Range("B6").Value2 = 0.05485232067510549
$ws.Range("D6").Value2 = 0.01687763713080169
$ws.Range("F6").Value2 = 0.0759493670886076
$ws.Range("J6").Value2 = 0.189873417721519
$ws.Range("O6").Value2 = 0.02531645569620253
$ws.Range("Q6").Value2 = 0.2151898734177215
$ws.Range("R6").Value2 = 0.06329113924050633
$ws.Range("S6").Value2 = 0.3586497890295359
$ws.Range("B7").Value2 = 0.1042944785276074
$ws.Range("D7").Value2 = 0.0245398773006135
$ws.Range("F7").Value2 = 0.06134969325153374
$ws.Range("J7").Value2 = 0.1042944785276074
$ws.Range("O7").Value2 = 0.03680981595092025
$ws.Range("Q7").Value2 = 0.1901840490797546
$ws.Range("R7").Value2 = 0.07975460122699386
$ws.Range("S7").Value2 = 0.3987730061349693
$ws.Range("B8").Value2 = 0.09018036072144289
$ws.Range("D8").Value2 = 0.01603206412825651
$ws.Range("F8").Value2 = 0.05210420841683366
$ws.Range("J8").Value2 = 0.1222444889779559
$ws.Range("O8").Value2 = 0.01603206412825651
$ws.Range("Q8").Value2 = 0.2044088176352706
$ws.Range("R8").Value2 = 0.09819639278557114
$ws.Range("S8").Value2 = 0.4008016032064128
$ws.Range("B9").Value2 = 0.1272727272727273
$ws.Range("D9").Value2 = 0.02272727272727273
$ws.Range("E9").Value2 = 0.004545454545454545
$ws.Range("F9").Value2 = 0.07727272727272727
$ws.Range("J9").Value2 = 0.1227272727272727
$ws.Range("O9").Value2 = 0.004545454545454545
$ws.Range("Q9").Value2 = 0.2045454545454546
$ws.Range("R9").Value2 = 0.06818181818181818
$ws.Range("S9").Value2 = 0.3681818181818182
$ws.Range("B10").Value2 = 0.1020819341840161
$ws.Range("D10").Value2 = 0.02014775016789792
$ws.Range("E10").Value2 = 0.000671591672263264
$ws.Range("F10").Value2 = 0.06984553391537945
$ws.Range("J10").Value2 = 0.1155137676292814
$ws.Range("O10").Value2 = 0.02350570852921424
$ws.Range("Q10").Value2 = 0.2316991269308261
$ws.Range("R10").Value2 = 0.08596373404969779
$ws.Range("S10").Value2 = 0.3505708529214238
$ws.Range("F11").Value2 = 0.003891050583657588
$ws.Range("G11").Value2 = 0.1400778210116732
$ws.Range("J11").Value2 = 0.08171206225680934
$ws.Range("K11").Value2 = 0.178988326848249
$ws.Range("L11").Value2 = 0.5719844357976653
$ws.Range("S11").Value2 = 0.02334630350194553
$ws.Range("G12").Value2 = 0.726027397260274
$ws.Range("J12").Value2 = 0.2191780821917808
$ws.Range("K12").Value2 = 0.0136986301369863
$ws.Range("L12").Value2 = 0.0273972602739726
$ws.Range("S12").Value2 = 0.0136986301369863
$ws.Range("G13").Value2 = 0.6666666666666666
$ws.Range("J13").Value2 = 0.2777777777777778
$ws.Range("S13").Value2 = 0.05555555555555555
$ws.Range("F15").Value2 = 0.01886792452830189
$ws.Range("H15").Value2 = 0.1660377358490566
$ws.Range("I15").Value2 = 0.0339622641509434
$ws.Range("J15").Value2 = 0.4037735849056604
$ws.Range("K15").Value2 = 0.0339622641509434
$ws.Range("M15").Value2 = 0.007547169811320755
$ws.Range("O15").Value2 = 0.06037735849056604
$ws.Range("S15").Value2 = 0.2754716981132075
$ws.Range("F16").Value2 = 0.02884615384615385
$ws.Range("H16").Value2 = 0.1875
$ws.Range("I16").Value2 = 0.05288461538461538
$ws.Range("J16").Value2 = 0.4423076923076923
$ws.Range("K16").Value2 = 0.0673076923076923
$ws.Range("M16").Value2 = 0.01923076923076923
$ws.Range("O16").Value2 = 0.07692307692307693
$ws.Range("S16").Value2 = 0.125
$ws.Range("F17").Value2 = 0.008771929824561403
$ws.Range("H17").Value2 = 0.1771929824561403
$ws.Range("I17").Value2 = 0.1052631578947368
$ws.Range("J17").Value2 = 0.4403508771929824
$ws.Range("K17").Value2 = 0.08947368421052632
$ws.Range("M17").Value2 = 0.02105263157894737
$ws.Range("O17").Value2 = 0.06842105263157895
$ws.Range("S17").Value2 = 0.08947368421052632
$ws.Range("F18").Value2 = 0.0045662100456621
$ws.Range("H18").Value2 = 0.1598173515981735
$ws.Range("I18").Value2 = 0.0639269406392694
$ws.Range("J18").Value2 = 0.5205479452054794
$ws.Range("K18").Value2 = 0.0639269406392694
$ws.Range("M18").Value2 = 0.0228310502283105
$ws.Range("N18").Value2 = 0.0045662100456621
$ws.Range("O18").Value2 = 0.0776255707762557
$ws.Range("S18").Value2 = 0.0821917808219178
$ws.Range("F19").Value2 = 0.00881057268722467
$ws.Range("H19").Value2 = 0.2063142437591777
$ws.Range("I19").Value2 = 0.09251101321585903
$ws.Range("J19").Value2 = 0.3986784140969163
$ws.Range("K19").Value2 = 0.08883994126284875
$ws.Range("M19").Value2 = 0.01101321585903084
$ws.Range("N19").Value2 = 0.001468428781204112
$ws.Range("O19").Value2 = 0.06607929515418502
$ws.Range("S19").Value2 = 0.1262848751835536
